$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (row 1 is the header). This shifts every
# existing data row (old row 2 -> new row 3, old row 44 -> new row 45) down
# by one, matching the dimension change from A1:R44 to A1:R45.
$ws.Rows(2).Insert()

# The freshly inserted row inherited bold/centered formatting from the
# header row above it; strip that so it matches the plain formatting used
# by the rest of the data rows.
$ws.Range("A2:R2").ClearFormats()

# Re-apply the date display format used by the other rows' "Fecha" column.
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

# Populate the new week's record (new latest observation for this
# market/product combination).
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44882
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 100112026
$ws.Cells.Item(2, 7).Value = "Haba"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 100
$ws.Cells.Item(2, 11).Value = 7500
$ws.Cells.Item(2, 12).Value = 8000
$ws.Cells.Item(2, 13).Value = 7750
$ws.Cells.Item(2, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(2, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 16).Value = 310
$ws.Cells.Item(2, 17).Value = 25
$ws.Cells.Item(2, 18).Value = "Hortaliza"
